$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 796.75
$ws.Range("I6").Value = 594.5
$ws.Range("K6").Value = 1783.5
$ws.Range("M6").Value = -1671.5
$ws.Range("H33").Value = 153.5625
$ws.Range("I33").Value = 154.07143
$ws.Range("J33").Value = 150
$ws.Range("K33").Value = 154.07143
$ws.Range("L33").Value = 150
$ws.Range("M33").Value = 74.92857000000001
$ws.Range("N33").Value = -608
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H64").Value = 4498.25
$ws.Range("J64").Value = 4831.3335
$ws.Range("L64").Value = 4831.3335
$ws.Range("N64").Value = -5327.3335
$ws.Range("H67").Value = 4498.25
$ws.Range("J67").Value = 4831.3335
$ws.Range("L67").Value = 4831.3335
$ws.Range("N67").Value = -6547.3335
$ws.Range("H107").Value = 420.58334
$ws.Range("I107").Value = 416.2
$ws.Range("K107").Value = 416.2
$ws.Range("M107").Value = 1503.8
$ws.Range("H112").Value = 4123.909
$ws.Range("J112").Value = 4176.3
$ws.Range("L112").Value = 12528.9
$ws.Range("N112").Value = -14744.9
$ws.Range("H116").Value = 4949.375
$ws.Range("I116").Value = 3866
$ws.Range("J116").Value = 5599.4
$ws.Range("K116").Value = 3866
$ws.Range("L116").Value = 5599.4
$ws.Range("M116").Value = -424
$ws.Range("N116").Value = -12483.4
$ws.Range("H132").Value = 3242.3572
$ws.Range("I132").Value = 3261
$ws.Range("K132").Value = 9783
$ws.Range("M132").Value = -7253

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2591.2
$ws.Range("I88").Value = 3506
$ws.Range("J88").Value = 2362.5
$ws.Range("K88").Value = 3506
$ws.Range("L88").Value = 2362.5
$ws.Range("M88").Value = -3100
$ws.Range("N88").Value = -3174.5
$ws.Range("H91").Value = 2591.2
$ws.Range("I91").Value = 3506
$ws.Range("J91").Value = 2362.5
$ws.Range("K91").Value = 3506
$ws.Range("L91").Value = 2362.5
$ws.Range("M91").Value = -2102
$ws.Range("N91").Value = -5170.5
$ws.Range("H102").Value = 4975
$ws.Range("I102").Value = 4950
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 4950
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -3328
$ws.Range("N102").Value = -8244
$ws.Range("H110").Value = 4112063.2
$ws.Range("I110").Value = 5286569
$ws.Range("J110").Value = 1294
$ws.Range("K110").Value = 5286569
$ws.Range("L110").Value = 1294
$ws.Range("M110").Value = -5284524
$ws.Range("N110").Value = -5384
$ws.Range("H132").Value = 1666.3334
$ws.Range("I132").Value = 1500
$ws.Range("K132").Value = 4500
$ws.Range("M132").Value = -1970

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2356.1428
$ws.Range("I20").Value = 2356.1428
$ws.Range("K20").Value = 2356.1428
$ws.Range("M20").Value = -2109.1428
$ws.Range("H86").Value = 1301.625
$ws.Range("I86").Value = 1262.1666
$ws.Range("K86").Value = 1262.1666
$ws.Range("M86").Value = -139.1666
$ws.Range("H89").Value = 1301.625
$ws.Range("I89").Value = 1262.1666
$ws.Range("K89").Value = 6310.833000000001
$ws.Range("M89").Value = -694.8330000000005
$ws.Range("H107").Value = 1912.125
$ws.Range("J107").Value = 1497.6666
$ws.Range("L107").Value = 1497.6666
$ws.Range("N107").Value = -5337.6666
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 79000
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 79000
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H86").Value = 14665.333
$ws.Range("I86").Value = 14665.333
$ws.Range("K86").Value = 14665.333
$ws.Range("M86").Value = -13542.333
$ws.Range("H89").Value = 14665.333
$ws.Range("I89").Value = 14665.333
$ws.Range("K89").Value = 73326.66500000001
$ws.Range("M89").Value = -67710.66500000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 16737.375
$ws.Range("I106").Value = 11300
$ws.Range("K106").Value = 33900
$ws.Range("M106").Value = -32954
$ws.Range("H113").Value = 1404.25
$ws.Range("I113").Value = 1314
$ws.Range("K113").Value = 3942
$ws.Range("M113").Value = -1772
$ws.Range("H121").Value = 16550.375
$ws.Range("I121").Value = 27567
$ws.Range("K121").Value = 82701
$ws.Range("M121").Value = -81391
$ws.Range("H122").Value = 896.6667
$ws.Range("I122").Value = 895
$ws.Range("J122").Value = 897.5
$ws.Range("K122").Value = 8055
$ws.Range("L122").Value = 8077.5
$ws.Range("M122").Value = -5605
$ws.Range("N122").Value = -12977.5
$ws.Range("H131").Value = 771522.1
$ws.Range("I131").Value = 2058.6
$ws.Range("J131").Value = 1252436.9
$ws.Range("K131").Value = 6175.799999999999
$ws.Range("L131").Value = 3757310.7
$ws.Range("M131").Value = -1135.799999999999
$ws.Range("N131").Value = -3767390.7
$ws.Range("H132").Value = 1900
$ws.Range("I132").Value = 1900
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 17100
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -14570
$ws.Range("N132").ClearContents()
$ws.Range("H133").Value = 22861.857
$ws.Range("I133").Value = 15000
$ws.Range("J133").Value = 24172.166
$ws.Range("K133").Value = 45000
$ws.Range("L133").Value = 72516.49800000001
$ws.Range("M133").Value = -39940
$ws.Range("N133").Value = -82636.49800000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 45760
$ws.Range("J141").Value = 45760
$ws.Range("L141").Value = 45760
$ws.Range("N141").Value = -56120

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2416.1667
$ws.Range("J46").Value = 2377.5
$ws.Range("L46").Value = 2377.5
$ws.Range("N46").Value = -2753.5
$ws.Range("H82").Value = 2000
$ws.Range("I82").Value = 2000
$ws.Range("K82").Value = 2000
$ws.Range("M82").Value = -1639
$ws.Range("H85").Value = 2000
$ws.Range("I85").Value = 2000
$ws.Range("K85").Value = 2000
$ws.Range("M85").Value = -752

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 22650
$ws.Range("J62").Value = 12299.667
$ws.Range("L62").Value = 12299.667
$ws.Range("N62").Value = -13547.667
$ws.Range("H65").Value = 22650
$ws.Range("J65").Value = 12299.667
$ws.Range("L65").Value = 61498.335
$ws.Range("N65").Value = -67738.33499999999
$ws.Range("H88").Value = 13997
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 13997
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 13997
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -14809
$ws.Range("H91").Value = 13997
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 13997
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 13997
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -16805
$ws.Range("H100").Value = 10000273
$ws.Range("I100").Value = 12500286
$ws.Range("J100").Value = 221
$ws.Range("K100").Value = 25000572
$ws.Range("L100").Value = 442
$ws.Range("M100").Value = -25000031
$ws.Range("N100").Value = -1524
